# Auto-generated edit script: updates Leve profit-tracking cells
# per the commit diff (scheduled market-price data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 1801
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1801
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 1801
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -2151
# Row 64 (Leve Item ID 5506)
$ws.Range("H64").Value = 3095.0625
$ws.Range("I64").Value = 2922.76
$ws.Range("J64").Value = 3710.4285
$ws.Range("K64").Value = 2922.76
$ws.Range("L64").Value = 3710.4285
$ws.Range("M64").Value = -2674.76
$ws.Range("N64").Value = -4206.4285
# Row 67 (Leve Item ID 5506)
$ws.Range("H67").Value = 3095.0625
$ws.Range("I67").Value = 2922.76
$ws.Range("J67").Value = 3710.4285
$ws.Range("K67").Value = 2922.76
$ws.Range("L67").Value = 3710.4285
$ws.Range("M67").Value = -2064.76
$ws.Range("N67").Value = -5426.4285
# Row 69 (Leve Item ID 12616)
$ws.Range("H69").Value = 2842.8572
$ws.Range("I69").Value = 1850
$ws.Range("J69").Value = 4166.6665
$ws.Range("K69").Value = 5550
$ws.Range("L69").Value = 12499.9995
$ws.Range("M69").Value = -4676
$ws.Range("N69").Value = -14247.9995
# Row 72 (Leve Item ID 12616)
$ws.Range("H72").Value = 2842.8572
$ws.Range("I72").Value = 1850
$ws.Range("J72").Value = 4166.6665
$ws.Range("K72").Value = 16650
$ws.Range("L72").Value = 37499.9985
$ws.Range("M72").Value = -12282
$ws.Range("N72").Value = -46235.9985
# Row 96 (Leve Item ID 19894)
$ws.Range("H96").Value = 450.2857
$ws.Range("I96").Value = 190.4
$ws.Range("J96").Value = 1100
$ws.Range("K96").Value = 571.2
$ws.Range("L96").Value = 3300
$ws.Range("M96").Value = 801.8
$ws.Range("N96").Value = -6046
# Row 107 (Leve Item ID 27766)
$ws.Range("H107").Value = 1043.5
$ws.Range("I107").Value = 1219
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 1219
$ws.Range("L107").Value = 400
$ws.Range("M107").Value = 701
$ws.Range("N107").Value = -4240
# Row 115 (Leve Item ID 27957)
$ws.Range("H115").Value = 822.7778
$ws.Range("I115").Value = 300.85715
$ws.Range("J115").Value = 1154.909
$ws.Range("K115").Value = 902.5714499999999
$ws.Range("L115").Value = 3464.727
$ws.Range("M115").Value = 664.4285500000001
$ws.Range("N115").Value = -6598.727000000001
# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 37860.594
$ws.Range("I132").Value = 40836.28
$ws.Range("J132").Value = 664.5
$ws.Range("K132").Value = 122508.84
$ws.Range("L132").Value = 1993.5
$ws.Range("M132").Value = -119978.84
$ws.Range("N132").Value = -7053.5

$ws = $wb.Worksheets.Item("ARM")
# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 25154.465
$ws.Range("I122").Value = 1916.9333
$ws.Range("J122").Value = 78779.53999999999
$ws.Range("K122").Value = 5750.7999
$ws.Range("L122").Value = 236338.62
$ws.Range("M122").Value = -3300.7999
$ws.Range("N122").Value = -241238.62

$ws = $wb.Worksheets.Item("BSM")
# Row 64 (Leve Item ID 14184)
$ws.Range("H64").Value = 429.5263
$ws.Range("I64").Value = 467
$ws.Range("J64").Value = 412.23077
$ws.Range("K64").Value = 467
$ws.Range("L64").Value = 412.23077
$ws.Range("M64").Value = -242
$ws.Range("N64").Value = -862.23077
# Row 67 (Leve Item ID 14184)
$ws.Range("H67").Value = 429.5263
$ws.Range("I67").Value = 467
$ws.Range("J67").Value = 412.23077
$ws.Range("K67").Value = 467
$ws.Range("L67").Value = 412.23077
$ws.Range("M67").Value = 313
$ws.Range("N67").Value = -1972.23077
# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 50052456
$ws.Range("I86").Value = 83335850
$ws.Range("J86").Value = 127366.875
$ws.Range("K86").Value = 83335850
$ws.Range("L86").Value = 127366.875
$ws.Range("M86").Value = -83334727
$ws.Range("N86").Value = -129612.875
# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 50052456
$ws.Range("I89").Value = 83335850
$ws.Range("J89").Value = 127366.875
$ws.Range("K89").Value = 416679250
$ws.Range("L89").Value = 636834.375
$ws.Range("M89").Value = -416673634
$ws.Range("N89").Value = -648066.375

$ws = $wb.Worksheets.Item("CUL")
# Row 107 (Leve Item ID 27838)
$ws.Range("H107").Value = 234.81818
$ws.Range("I107").Value = 242
$ws.Range("J107").Value = 232.97144
$ws.Range("K107").Value = 726
$ws.Range("L107").Value = 698.91432
$ws.Range("M107").Value = 1194
$ws.Range("N107").Value = -4538.91432

$ws = $wb.Worksheets.Item("GSM")
# Row 5 (Leve Item ID 1681)
$ws.Range("H5").Value = 11970.5
$ws.Range("I5").Value = 10000
$ws.Range("J5").Value = 12189.444
$ws.Range("K5").Value = 10000
$ws.Range("L5").Value = 12189.444
$ws.Range("M5").Value = -9888
$ws.Range("N5").Value = -12413.444
# Row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 3048
$ws.Range("I102").Value = 2072.5
$ws.Range("J102").Value = 4999
$ws.Range("K102").Value = 2072.5
$ws.Range("L102").Value = 4999
$ws.Range("M102").Value = -450.5
$ws.Range("N102").Value = -8243
# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 1340.6923
$ws.Range("I122").Value = 1322.9333
$ws.Range("J122").Value = 1364.909
$ws.Range("K122").Value = 3968.7999
$ws.Range("L122").Value = 4094.727
$ws.Range("M122").Value = -1518.7999
$ws.Range("N122").Value = -8994.727000000001
# Row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 12829795
$ws.Range("I126").Value = 26498
$ws.Range("J126").Value = 18520150
$ws.Range("K126").Value = 79494
$ws.Range("L126").Value = 55560450
$ws.Range("M126").Value = -77024
$ws.Range("N126").Value = -55565390

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 2737.5557
$ws.Range("I7").Value = 1167.1666
$ws.Range("J7").Value = 3522.75
$ws.Range("K7").Value = 1167.1666
$ws.Range("L7").Value = 3522.75
$ws.Range("M7").Value = -1055.1666
$ws.Range("N7").Value = -3746.75
# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 2308.423
$ws.Range("I40").Value = 2144.3333
$ws.Range("J40").Value = 2677.625
$ws.Range("K40").Value = 2144.3333
$ws.Range("L40").Value = 2677.625
$ws.Range("M40").Value = -2008.3333
$ws.Range("N40").Value = -2949.625
# Row 81 (Leve Item ID 10897)
$ws.Range("H81").Value = 21000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 21000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 21000
$ws.Range("N81").Value = -22996
# Row 84 (Leve Item ID 10897)
$ws.Range("H84").Value = 21000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 21000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 63000
$ws.Range("N84").Value = -72984
# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 2737.5557
$ws.Range("I126").Value = 1167.1666
$ws.Range("J126").Value = 3522.75
$ws.Range("K126").Value = 3501.4998
$ws.Range("L126").Value = 10568.25
$ws.Range("M126").Value = -1031.4998
$ws.Range("N126").Value = -15508.25
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 29620.777
$ws.Range("I132").Value = 57565.332
$ws.Range("J132").Value = 1676.2222
$ws.Range("K132").Value = 172695.996
$ws.Range("L132").Value = 5028.6666
$ws.Range("M132").Value = -170165.996
$ws.Range("N132").Value = -10088.6666

$ws = $wb.Worksheets.Item("WVR")
# Row 8 (Leve Item ID 2999)
$ws.Range("H8").Value = 2003.6666
$ws.Range("I8").Value = 2003
$ws.Range("J8").Value = 2004
$ws.Range("K8").Value = 2003
$ws.Range("L8").Value = 2004
$ws.Range("M8").Value = -1863
$ws.Range("N8").Value = -2284
# Row 39 (Leve Item ID 3106)
$ws.Range("H39").Value = 5000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 5000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 5000
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -5826
# Row 80 (Leve Item ID 10911)
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
# Row 81 (Leve Item ID 12596)
$ws.Range("H81").Value = 3370.4
$ws.Range("I81").Value = 2675.8823
$ws.Range("J81").Value = 4846.25
$ws.Range("K81").Value = 5351.7646
$ws.Range("L81").Value = 9692.5
$ws.Range("M81").Value = -4290.7646
$ws.Range("N81").Value = -11814.5
# Row 83 (Leve Item ID 10911)
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
# Row 84 (Leve Item ID 12596)
$ws.Range("H84").Value = 3370.4
$ws.Range("I84").Value = 2675.8823
$ws.Range("J84").Value = 4846.25
$ws.Range("K84").Value = 26758.823
$ws.Range("L84").Value = 48462.5
$ws.Range("M84").Value = -21454.823
$ws.Range("N84").Value = -59070.5

